$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update TPM-derived metrics (ligand/receptor expression, specificity, and edge weights)
# to reflect the new TPM values, per row.

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3915976666666667
$ws.Range("H2").Value = 1.174793
$ws.Range("I2").Value = 0.02606065131430495
$ws.Range("J2").Value = 0.02606065131430495
$ws.Range("M2").Value = 29.52617166666667
$ws.Range("N2").Value = 88.57851500000001
$ws.Range("O2").Value = 0.3218391660320701
$ws.Range("P2").Value = 0.3218391660320701
$ws.Range("Q2").Value = 11.56237993026611
$ws.Range("R2").Value = 104.061419372395
$ws.Range("S2").Value = 0.008387338285248477
$ws.Range("T2").Value = 0.008387338285248477

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3915976666666667
$ws.Range("H3").Value = 1.174793
$ws.Range("I3").Value = 0.02606065131430495
$ws.Range("J3").Value = 0.02606065131430495
$ws.Range("O3").Value = 0.4328989896002822
$ws.Range("P3").Value = 0.4328989896002822
$ws.Range("Q3").Value = 15.552310338413
$ws.Range("R3").Value = 139.970793045717
$ws.Range("S3").Value = 0.01128162962228788
$ws.Range("T3").Value = 0.01128162962228788

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3915976666666667
$ws.Range("H4").Value = 1.174793
$ws.Range("I4").Value = 0.02606065131430495
$ws.Range("J4").Value = 0.02606065131430495
$ws.Range("M4").Value = 22.50081433333333
$ws.Range("N4").Value = 67.502443
$ws.Range("O4").Value = 0.2452618443676477
$ws.Range("P4").Value = 0.2452618443676476
$ws.Range("Q4").Value = 8.811266391033223
$ws.Range("R4").Value = 79.30139751929902
$ws.Range("S4").Value = 0.006391683406768594
$ws.Range("T4").Value = 0.006391683406768594

# Row 5
$ws.Range("I5").Value = 0.4187506438669658
$ws.Range("J5").Value = 0.4187506438669658
$ws.Range("M5").Value = 29.52617166666667
$ws.Range("N5").Value = 88.57851500000001
$ws.Range("O5").Value = 0.3218391660320701
$ws.Range("P5").Value = 0.3218391660320701
$ws.Range("Q5").Value = 185.7879138183984
$ws.Range("R5").Value = 1672.091224365585
$ws.Range("S5").Value = 0.1347703579975366
$ws.Range("T5").Value = 0.1347703579975366

# Row 6
$ws.Range("I6").Value = 0.4187506438669658
$ws.Range("J6").Value = 0.4187506438669658
$ws.Range("O6").Value = 0.4328989896002822
$ws.Range("P6").Value = 0.4328989896002822
$ws.Range("S6").Value = 0.1812767306244771
$ws.Range("T6").Value = 0.1812767306244771

# Row 7
$ws.Range("I7").Value = 0.4187506438669658
$ws.Range("J7").Value = 0.4187506438669658
$ws.Range("M7").Value = 22.50081433333333
$ws.Range("N7").Value = 67.502443
$ws.Range("O7").Value = 0.2452618443676477
$ws.Range("P7").Value = 0.2452618443676476
$ws.Range("Q7").Value = 141.5821665402197
$ws.Range("R7").Value = 1274.239498861977
$ws.Range("S7").Value = 0.102703555244952
$ws.Range("T7").Value = 0.102703555244952

# Row 8
$ws.Range("G8").Value = 8.342485333333334
$ws.Range("H8").Value = 25.027456
$ws.Range("I8").Value = 0.5551887048187292
$ws.Range("J8").Value = 0.5551887048187292
$ws.Range("M8").Value = 29.52617166666667
$ws.Range("N8").Value = 88.57851500000001
$ws.Range("O8").Value = 0.3218391660320701
$ws.Range("P8").Value = 0.3218391660320701
$ws.Range("Q8").Value = 246.3216540786489
$ws.Range("R8").Value = 2216.89488670784
$ws.Range("S8").Value = 0.1786814697492849
$ws.Range("T8").Value = 0.1786814697492849

# Row 9
$ws.Range("G9").Value = 8.342485333333334
$ws.Range("H9").Value = 25.027456
$ws.Range("I9").Value = 0.5551887048187292
$ws.Range("J9").Value = 0.5551887048187292
$ws.Range("O9").Value = 0.4328989896002822
$ws.Range("P9").Value = 0.4328989896002822
$ws.Range("Q9").Value = 331.3219968904961
$ws.Range("R9").Value = 2981.897972014464
$ws.Range("S9").Value = 0.2403406293535172
$ws.Range("T9").Value = 0.2403406293535172

# Row 10
$ws.Range("G10").Value = 8.342485333333334
$ws.Range("H10").Value = 25.027456
$ws.Range("I10").Value = 0.5551887048187292
$ws.Range("J10").Value = 0.5551887048187292
$ws.Range("M10").Value = 22.50081433333333
$ws.Range("N10").Value = 67.502443
$ws.Range("O10").Value = 0.2452618443676477
$ws.Range("P10").Value = 0.2452618443676476
$ws.Range("Q10").Value = 187.7127135638898
$ws.Range("R10").Value = 1689.414422075008
$ws.Range("S10").Value = 0.136166605715927
$ws.Range("T10").Value = 0.136166605715927
